$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (id=8, "Enchanting is key" quest) ---
# D9 (instructions): append a new paragraph about Skill Type Crafting requirement
$ws.Range("D9").Value = @'
<p>Today's lesson is about <a href="/information/enchanting" target="_blank">enchanting</a> and <a href="/information/disenchanting" target="_blank">disenchanting</a>. Just like crafting, enchanting is wildly important. As you may know, you cannot buy enchanted gear from the shop. Some players will sell it on the <a href="/information/market-board" target="_blank">market place</a>.</p><p>More importantly, you want a <a href="/information/currencies" target="_blank">currency</a> you can get from it: Gold Dust. This currency is used in a lot of quests we are going to start doing soon. So lets get a lot of it.</p><p><strong>Desktop/Mobile:</strong></p><p><strong>- </strong>To enchant, from the crafting drop down select enchant.</p><p>- Here you have an item to select, a suffix and a prefix to select. The cost of enchanting can get widely expensive so exploration is required. Since you need items: weapons, Armour, spells and Rings, you will also want to be crafting. You can also re-enchant enchanted items but that gets more costly.</p><p>You can also just go to the shop and buy multiple items to enchant, like broken daggers.</p><p>- Next, when you have enchanted, ideally 75 items, next go to your inventory</p><p>- From the actions on the inventory section of your character sheet, click Disenchant All.</p><p>- This will raise your disenchanting skill, and slowly over time your enchanting skill as well.</p><p>- Repeat till you meet the quest requirements.</p><p>The <strong>Skill Type Crafting to level 10</strong> basically means, level any rafting skill you want to level 10 or higher.</p>
'@

# E9 (required_level): 250 -> 160
$ws.Range("E9").Value = 160

# I9 (required_secondary_skill_level): 12 -> 5
$ws.Range("I9").Value = 5

# J9 (required_skill_type): new requirement, Crafting type id = 1
$ws.Range("J9").Value = 1

# K9 (required_skill_type_level): new requirement, level 10
$ws.Range("K9").Value = 10

# X9 (required_gold_dust): 1000 -> 500
$ws.Range("X9").Value = 500

# AE9 (required_stats): 300 -> 280
$ws.Range("AE9").Value = 280

# --- Row 10 (id=9, "Traveling To Dungeons" quest) ---
# C10 (intro_text / story): typo fixes
$ws.Range("C10").Value = @'
You have come a long way since you first set out. Learning more of the world, learning more about your self and developing your skills.<br /> <br /> You sit under the tree outside of town and watch the merchants go by. A few stop and ask if you would like to purchase some of their wares, a few have some interesting knickknacks, and a few others have some food worth buying.<br /> <br /> As the day passes by and the sun begins to set you think of heading back to town to wash, relax and get a good nights sleep.<br /> <br /> “Hello there” comes an unfamiliar voice.<br /> <br /> You look and see a man in a Fedora, Green Tunic and black leather pants. He seems older, but moves as if he is younger then his age.<br /> <br /> He comes closer and introduces himself, “I am The Poet and have been sent by The Guide to ask a favor of you, one from him.”<br /> <br /> You ask why he doesn’t come himself and ask this favor.<br /> <br /> “There is a darkness that is seeping from the depths into all the Planes. It is corrupting the Planes and causing the creatures to become vile and wicked. The Guide is busy investigating something deep with in Dungeons. Where he wants you to meet him.”<br /> <br /> You remember how you got to labyrinth, how hard can it be to get to Dungeons? You agree and The Poet tells you how to set off.
'@

# D10 (instructions): "Ring Lord" -> "Dark Enchantress"
$ws.Range("D10").Value = @'
<p>This is another <a href="/information/quests" target="_blank">quest</a> where we are going to a new <a href="/information/planes" target="_blank">plane</a>. How ever, to get there we have to complete a few quests.</p><p>Quests are how a lot of Tlessas features are unlocked for players. These quests, much like Guide Quests, will have a story – before and after completion – as well as explicit instructions on how to complete them, much like these.</p><p>Since you already know how to <a href="/information/traverse" target="_blank">traverse</a>, we just need to get you a quest item that allows you to traverse down to Dungeons.</p><p><strong>Note:</strong> The required secondary quest item, will require you to complete a quest from Labyrinth under the One Off Quest section called: <strong>Dark Enchantress.</strong> Click the quest, and click on the Requirements tab to see what you have to do.</p><p><strong>Desktop</strong></p><p>- Click the Quests tab</p><p>- Here you can see a list of quests for the plane you are on (if you are not on surface, select the Surface plane from the Planes drop down.)</p><p>- On the Surface quests you will see a quest (to the left) called Light The Way, complete the quests working down the tree until you complete it. Locked Quests are red, quests you can complete are Blue and Completed Quests are Green.</p><p>- Click on the quest to read the story, see the requirements and the reward.</p><p>- The Requirements tab will have explicit instructions beside each requirement.</p><p>- Much like Guide Quests, you can see your completed quests in the sidebar (Hamburger menu to the top left) under: Quest Log → Completed Quests.</p><p><strong>Mobile</strong></p><p>- Tap the Quests tab</p><p>- Here you can see a list of quests for the plane you are on (if you are not on surface, select the Surface plane from the Planes drop down.)</p><p>- On the Surface quests you will see a quest (to the left – you can scroll left/right) called Light The Way, complete the quests working down the tree until you complete it. Locked Quests are red, quests you can complete are Blue and Completed Quests are Green.</p><p>- Tap on the quest to read the story, see the requirements and the reward.</p><p>- The Requirements tab will have explicit instructions beside each requirement.</p><p>- Much like Guide Quests, you can see your completed quests in the sidebar (Hamburger menu to the top left) under: Quest Log → Completed Quests.</p><p>Quests, much like Enchanting are vitally important to progress not only the over all story of the game, but also to get access to new features, new planes and so on of the game.</p><p>Each Plane will tell a story that all come together to tell a general story of The Creator (not the game Creator – an NPC) Quest chains going straight down tell their own story.</p>
'@

# E10 (required_level): 300 -> 200
$ws.Range("E10").Value = 200

# H10 (required_secondary_skill "Ring Crafting"): removed entirely
$ws.Range("H10").ClearContents()

# I10 (required_secondary_skill_level): 10 -> 0
$ws.Range("I10").Value = 0

# S10 (required_game_map_id): new requirement "Dungeons"
$ws.Range("S10").Value = "Dungeons"

# T10 (required_quest_id): new requirement "Light the way"
$ws.Range("T10").Value = "Light the way"

# U10 (required_quest_item_id): stays "Torch" (unchanged text, re-set for safety)
$ws.Range("U10").Value = "Torch"

# V10 (secondary_quest_item_id): "Ring Crafter's Master Book" -> "Enchantress' Diary"
$ws.Range("V10").Value = "Enchantress' Diary"

# AE10 (required_stats): 400 -> 500
$ws.Range("AE10").Value = 500
